# ---------------------------------------------------------------------------
# This script applies the "ADDITIONAL SCRAPING" edit to the workbook:
#   1. Adds a new "Player Info" sheet (at the very front) with player bio data
#   2. Keeps "ODI Batting" (renames its MATCH_CARD_LINK column to MATCH_CODE,
#      converting the stored URLs into bare match-code numbers, and drops a
#      handful of stray empty INNING_NUMBER cells on "did not bat" rows)
#   3. Keeps "ODI Bowling" (same MATCH_CARD_LINK -> MATCH_CODE treatment)
#   4. Adds a new "ODI Batting Extra" sheet (at the very end) with extra
#      per-match batting stats
# Final sheet order: Player Info, ODI Batting, ODI Bowling, ODI Batting Extra
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

$wsBatting = $wb.Worksheets.Item("ODI Batting")
$wsBowling = $wb.Worksheets.Item("ODI Bowling")

# ---------------------------------------------------------------------------
# Helper: apply the same bold / bordered / centered-top header look used by
# the header row of the existing sheets to a given cell.
# ---------------------------------------------------------------------------
function Set-HeaderStyle($cell) {
    $cell.Font.Bold = $true
    $cell.Borders.LineStyle = 1
    $cell.HorizontalAlignment = -4108
    $cell.VerticalAlignment = -4160
}

# ---------------------------------------------------------------------------
# 1) "ODI Batting": MATCH_CARD_LINK -> MATCH_CODE
# ---------------------------------------------------------------------------
$battingLastRow = $wsBatting.UsedRange.Rows.Count
for ($r = 2; $r -le $battingLastRow; $r++) {
    $cell = $wsBatting.Cells.Item($r, 4)
    $url = $cell.Value2
    $parts = $url -split "MatchCode="
    $code = $parts[1]
    $cell.NumberFormat = "@"
    $cell.Value = $code
}
$wsBatting.Cells.Item(1, 4).Value = "MATCH_CODE"

# Rows where the batsman "did not bat" used to carry a stray empty
# INNING_NUMBER (column B) cell; remove those leftover blank cells.
$emptyInningRows = @(5, 23, 65, 91, 115, 126, 140)
foreach ($r in $emptyInningRows) {
    $wsBatting.Cells.Item($r, 2).ClearContents()
}

# ---------------------------------------------------------------------------
# 2) "ODI Bowling": MATCH_CARD_LINK -> MATCH_CODE
# ---------------------------------------------------------------------------
$bowlingLastRow = $wsBowling.UsedRange.Rows.Count
for ($r = 2; $r -le $bowlingLastRow; $r++) {
    $cell = $wsBowling.Cells.Item($r, 2)
    $url = $cell.Value2
    $parts = $url -split "MatchCode="
    $code = $parts[1]
    $cell.NumberFormat = "@"
    $cell.Value = $code
}
$wsBowling.Cells.Item(1, 2).Value = "MATCH_CODE"

# ---------------------------------------------------------------------------
# 3) New "Player Info" sheet, inserted before "ODI Batting"
# ---------------------------------------------------------------------------
$wsPlayerInfo = $wb.Worksheets.Add($wsBatting)
$wsPlayerInfo.Name = "Player Info"

$playerInfoHeaders = @("ID", "NAME", "BATTING_HAND", "BOWL_STYLE")
for ($c = 1; $c -le $playerInfoHeaders.Length; $c++) {
    $cell = $wsPlayerInfo.Cells.Item(1, $c)
    $cell.Value = $playerInfoHeaders[$c - 1]
    Set-HeaderStyle $cell
}

$wsPlayerInfo.Cells.Item(2, 1).NumberFormat = "@"
$wsPlayerInfo.Cells.Item(2, 1).Value = "3847"
$wsPlayerInfo.Cells.Item(2, 2).Value = "Francois Du Plessis"
$wsPlayerInfo.Cells.Item(2, 3).Value = "Right Handed"
$wsPlayerInfo.Cells.Item(2, 4).Value = "Right Arm Leg Break"

# ---------------------------------------------------------------------------
# 4) New "ODI Batting Extra" sheet, inserted after "ODI Bowling" (appended
#    as the last sheet)
# ---------------------------------------------------------------------------
$wsBowlingFresh = $wb.Worksheets.Item("ODI Bowling")
$wsExtra = $wb.Worksheets.Add($null, $wsBowlingFresh)
$wsExtra.Name = "ODI Batting Extra"

$extraHeaders = @("MATCH_CODE", "BATTING_POSITION", "NUM_4", "NUM_6", "PERCENT_RUNS_OF_TOTAL", "MAN_OF_MATCH")
for ($c = 1; $c -le $extraHeaders.Length; $c++) {
    $cell = $wsExtra.Cells.Item(1, $c)
    $cell.Value = $extraHeaders[$c - 1]
    Set-HeaderStyle $cell
}

$extraData = @(
    ,@("4226", "4", "15", "2", "39.06%", "NO")
    ,@("4237", "5", $null, $null, $null, "NO")
    ,@("4238", "4", "2", "0", "3.86%", "NO")
    ,@("4241", "4", "3", "1", "21.39%", "NO")
    ,@("4244", "4", "5", "1", "34.76%", "NO")
    ,@("4247", $null, $null, $null, $null, "NO")
    ,@("4261", "3", "15", "1", "48.28%", "YES")
    ,@("4264", "3", "7", "0", "22.71%", "NO")
    ,@("4269", "3", "7", "0", "10.88%", "NO")
    ,@("4271", $null, $null, $null, $null, "NO")
    ,@("4272", "3", "1", "0", "17.78%", "NO")
    ,@("4303", $null, $null, $null, $null, "NO")
    ,@("4307", $null, $null, $null, $null, "NO")
    ,@("4310", $null, $null, $null, $null, "NO")
    ,@("4317", $null, $null, $null, $null, "NO")
    ,@("4323", "5", $null, $null, $null, "NO")
    ,@("4328", "3", "4", "0", "9.54%", "NO")
    ,@("4334", "3", "5", "0", "24.32%", "NO")
    ,@("4339", "3", "10", "1", "46.60%", "NO")
    ,@("4351", "3", "7", "2", "30.77%", "YES")
)

for ($i = 0; $i -lt $extraData.Length; $i++) {
    $row = $extraData[$i]
    $r = $i + 2

    $matchCodeCell = $wsExtra.Cells.Item($r, 1)
    $matchCodeCell.NumberFormat = "@"
    $matchCodeCell.Value = $row[0]

    $posCell = $wsExtra.Cells.Item($r, 2)
    if ($row[1] -eq $null) {
        $posCell.Value = ""
    } else {
        $posCell.Value = [int]$row[1]
    }

    $num4Cell = $wsExtra.Cells.Item($r, 3)
    $num4Cell.NumberFormat = "@"
    if ($row[2] -eq $null) {
        $num4Cell.Value = ""
    } else {
        $num4Cell.Value = $row[2]
    }

    $num6Cell = $wsExtra.Cells.Item($r, 4)
    $num6Cell.NumberFormat = "@"
    if ($row[3] -eq $null) {
        $num6Cell.Value = ""
    } else {
        $num6Cell.Value = $row[3]
    }

    $pctCell = $wsExtra.Cells.Item($r, 5)
    $pctCell.NumberFormat = "@"
    if ($row[4] -eq $null) {
        $pctCell.Value = ""
    } else {
        $pctCell.Value = $row[4]
    }

    $momCell = $wsExtra.Cells.Item($r, 6)
    $momCell.Value = $row[5]
}
